# Week 6 timesheet update (attendance up to date as of 10/04/2021).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read/Study (row 7): logged 2 hours on Saturday (H)
$ws.Range("H7").Value = 2

# Team Work (row 11): logged 1 hour on Friday (F)
$ws.Range("F11").Value = 1

# Organizing (misc.) (row 12): logged 1 hour on Saturday (H)
$ws.Range("H12").Value = 1

# Leave the selection where the author last left it
$ws.Range("L15").Select()
